$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# cells on the zh-cn sheet for the 3f94c6a7... handback row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-11 16:29:21"
$wsZhCn.Range("H3").Value = "2016-03-11 16:29:40"

# Update the same cells on the de-de sheet for the 3f94c6a7... handback row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-11 16:29:24"
$wsDeDe.Range("H3").Value = "2016-03-11 16:29:45"
